# Adds 2024-12-30 data to the violent crime workbook.
# For each affected sheet, update the specific cells whose values
# changed due to the additional day of data (and a couple of minor
# historical corrections to 2022 figures).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('K2').Value = 7867
$ws.Range('I3').Value = 7493
$ws.Range('K3').Value = 8152
$ws.Range('K4').Value = 1716
$ws.Range('K6').Value = 9082
$ws.Range('I7').Value = 26279
$ws.Range('K7').Value = 27395

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('K3').Value = 543
$ws.Range('K6').Value = 603
$ws.Range('K7').Value = 1790

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('K2').Value = 203
$ws.Range('K4').Value = 26
$ws.Range('K7').Value = 584

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('K3').Value = 409
$ws.Range('K6').Value = 363
$ws.Range('K7').Value = 1148

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('K2').Value = 154
$ws.Range('K7').Value = 455

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('K3').Value = 300
$ws.Range('K6').Value = 272
$ws.Range('K7').Value = 903

$ws = $wb.Worksheets.Item('New City')
$ws.Range('K2').Value = 207
$ws.Range('K3').Value = 151
$ws.Range('K7').Value = 636

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('K6').Value = 113
$ws.Range('K7').Value = 465

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range('K6').Value = 44
$ws.Range('K7').Value = 112

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('K2').Value = 238
$ws.Range('K7').Value = 815
$ws.Range('K8').Value = 1790
$ws.Range('K14').Value = 129
$ws.Range('K15').Value = 278
$ws.Range('K16').Value = 68
$ws.Range('K18').Value = 184
$ws.Range('K19').Value = 788
$ws.Range('K20').Value = 672
$ws.Range('K26').Value = 35
$ws.Range('K29').Value = 1509
$ws.Range('K30').Value = 112
$ws.Range('K32').Value = 29
$ws.Range('K33').Value = 1148
$ws.Range('K34').Value = 156
$ws.Range('K37').Value = 903
$ws.Range('K42').Value = 1018
$ws.Range('K43').Value = 226
$ws.Range('K47').Value = 192
$ws.Range('K51').Value = 357
$ws.Range('K54').Value = 534
$ws.Range('K57').Value = 113
$ws.Range('I63').Value = 241
$ws.Range('K63').Value = 76
$ws.Range('K65').Value = 636
$ws.Range('K67').Value = 1066
$ws.Range('K68').Value = 71
$ws.Range('K71').Value = 83
$ws.Range('K75').Value = 91
$ws.Range('K76').Value = 379
$ws.Range('K79').Value = 665
$ws.Range('K83').Value = 584
$ws.Range('K85').Value = 1268
$ws.Range('K87').Value = 54
$ws.Range('K88').Value = 288
$ws.Range('K89').Value = 413
$ws.Range('K94').Value = 366
$ws.Range('K95').Value = 455
$ws.Range('K96').Value = 298
$ws.Range('K97').Value = 226
$ws.Range('K98').Value = 149
$ws.Range('K99').Value = 465
$ws.Range('I101').Value = 26279
$ws.Range('K101').Value = 27395

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('K6').Value = 301
$ws.Range('K7').Value = 1066

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('K2').Value = 84
$ws.Range('K7').Value = 534

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('K2').Value = 427
$ws.Range('K3').Value = 533
$ws.Range('K6').Value = 443
$ws.Range('K7').Value = 1509

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('K2').Value = 230
$ws.Range('K3').Value = 236
$ws.Range('K6').Value = 263
$ws.Range('K7').Value = 788

$ws = $wb.Worksheets.Item('River North')
$ws.Range('K3').Value = 76
$ws.Range('K6').Value = 186
$ws.Range('K7').Value = 379

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('K2').Value = 41
$ws.Range('K7').Value = 129

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('K5').Value = 17
$ws.Range('K6').Value = 389
$ws.Range('K7').Value = 1018

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('K3').Value = 64
$ws.Range('K7').Value = 298

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('K6').Value = 176
$ws.Range('K7').Value = 665

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('K2').Value = 227
$ws.Range('K7').Value = 672

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('K6').Value = 48
$ws.Range('K7').Value = 184

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('K2').Value = 266
$ws.Range('K6').Value = 231
$ws.Range('K7').Value = 815

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('K6').Value = 44
$ws.Range('K7').Value = 156

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('K6').Value = 172
$ws.Range('K7').Value = 366

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('K2').Value = 54
$ws.Range('K6').Value = 59
$ws.Range('K7').Value = 192

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('K3').Value = 69
$ws.Range('K7').Value = 278

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range('K6').Value = 83
$ws.Range('K7').Value = 149

$ws = $wb.Worksheets.Item('East Village')
$ws.Range('K6').Value = 25
$ws.Range('K7').Value = 35

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('K2').Value = 72
$ws.Range('K7').Value = 238

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('K6').Value = 127
$ws.Range('K7').Value = 226

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('K2').Value = 75
$ws.Range('K7').Value = 288

$ws = $wb.Worksheets.Item('Galewood')
$ws.Range('K2').Value = 10
$ws.Range('K7').Value = 29

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('K2').Value = 112
$ws.Range('K3').Value = 127
$ws.Range('K6').Value = 123
$ws.Range('K7').Value = 413

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range('K3').Value = 31
$ws.Range('K7').Value = 91

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('K2').Value = 95
$ws.Range('K3').Value = 99
$ws.Range('K6').Value = 116
$ws.Range('K7').Value = 357

$ws = $wb.Worksheets.Item('North Park')
$ws.Range('K3').Value = 16
$ws.Range('K7').Value = 71

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range('K6').Value = 48
$ws.Range('K7').Value = 113

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range('K6').Value = 81
$ws.Range('K7').Value = 226

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('K2').Value = 422
$ws.Range('K6').Value = 312
$ws.Range('K7').Value = 1268

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range('K3').Value = 31
$ws.Range('K7').Value = 83

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range('K6').Value = 26
$ws.Range('K7').Value = 54

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range('K4').Value = 6
$ws.Range('K6').Value = 40
$ws.Range('K7').Value = 68
